$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = "CONTRACTS,B2B"
$ws.Range("K4").Value = "B2G"
$ws.Range("K5").Value = "B2B, B2C,"

$ws.Range("K6").Select()
